# Added more performance results when GlusterFS was added to the test boxes
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("jon-akka01")

# --- New environment-note section (mirrors the earlier "Environment for ..." blocks) ---
# Reuse the bold "section header" formatting already used by A1/A17/A29 instead of
# minting a brand new style.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("A42").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("A42").Value = "Environment for jon-akka01 on sadbox "

$ws.Range("A43").Value = "Installed gluster file system"
$ws.Range("A43").NumberFormat = "0"

$ws.Range("A44").Value = "Snapshots being saved to a gluster volume mapped to 2 other machines"
$ws.Range("A44").NumberFormat = "0"

# --- New timing table (mirrors the LoadAsync()/SaveAsync() tables above) ---
$ws.Range("A46").Value = "LoadAsync()"
$ws.Range("H46").Value = "SaveAsync()"
# Column A carries a sheet-level default style (id 2); these header cells use the
# plain/default "Normal" style instead, same as the A10/A22/A35 headers above.
$ws.Range("A46").Style = "Normal"

$ws.Range("A47").Value = "Start"
$ws.Range("B47").Value = "Finished"
$ws.Range("C47").Value = "Elapsed"
$ws.Range("D47").Value = "CPU%"
$ws.Range("H47").Value = "Start"
$ws.Range("I47").Value = "Finished"
$ws.Range("J47").Value = "Elapsed"
$ws.Range("K47").Value = "CPU%"
$ws.Range("L47").Value = "Disk"
$ws.Range("A47").Style = "Normal"

# Row 48 - LoadAsync() timing/result row
$ws.Range("A48").Value = 0.38646990740740739
$ws.Range("B48").Value = 0.38988425925925929
$ws.Range("C48").Formula = "=B48-A48"
$ws.Range("D48").Value = "Not all threads used, bulk were at 20%"

# Row 48/49 - SaveAsync() timing/result row (H48/I48 blank, H49/I49 filled)
$ws.Range("H48").NumberFormat = "h:mm:ss"
$ws.Range("I48").NumberFormat = "h:mm:ss"
$ws.Range("J48").Formula = "=I48-H48"
$ws.Range("K48").Value = "All threads MAX 85%"
$ws.Range("L48").Value = "30MBps"

$ws.Range("C49").Formula = "=B49-A49"
$ws.Range("H49").Value = 0.23126157407407408
$ws.Range("I49").Value = 0.23283564814814817
$ws.Range("J49").Formula = "=I49-H49"

# Time-of-day number formatting for the new data rows (matches rows 12/24/37 above)
$ws.Range("A48:C49").NumberFormat = "h:mm:ss"
$ws.Range("H48:J49").NumberFormat = "h:mm:ss"

# --- View state: scroll position / active selection moved along with the new rows ---
$ws.Activate()
$ws.Range("G51").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1
